$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 44477
$ws.Cells.Item(2, 9).Value = 'Primera'
$ws.Cells.Item(2, 10).Value = 80
$ws.Cells.Item(2, 11).Value = 800
$ws.Cells.Item(2, 12).Value = 800
$ws.Cells.Item(2, 13).Value = 800
$ws.Cells.Item(2, 14).Value = '$/kilo (volumen en unidades)'
$ws.Cells.Item(2, 15).Value = 'Perú'
$ws.Cells.Item(2, 16).Value = 800
# Row 3
$ws.Cells.Item(3, 4).Value = 44305
$ws.Cells.Item(3, 10).Value = 100
$ws.Cells.Item(3, 11).Value = 2500
$ws.Cells.Item(3, 12).Value = 2500
$ws.Cells.Item(3, 13).Value = 2500
$ws.Cells.Item(3, 15).Value = 'Perú'
$ws.Cells.Item(3, 16).Value = 2500
# Row 4
$ws.Cells.Item(4, 4).Value = 44504
$ws.Cells.Item(4, 10).Value = 200
# Row 5
$ws.Cells.Item(5, 4).Value = 44497
$ws.Cells.Item(5, 8).Value = 'Sin especificar'
$ws.Cells.Item(5, 9).Value = 'Primera'
$ws.Cells.Item(5, 10).Value = 250
$ws.Cells.Item(5, 11).Value = 800
$ws.Cells.Item(5, 12).Value = 800
$ws.Cells.Item(5, 13).Value = 800
$ws.Cells.Item(5, 14).Value = '$/kilo (volumen en unidades)'
$ws.Cells.Item(5, 15).Value = 'Perú'
$ws.Cells.Item(5, 16).Value = 800
# Row 6
$ws.Cells.Item(6, 9).Value = 'Extra'
$ws.Cells.Item(6, 10).Value = 340
$ws.Cells.Item(6, 11).Value = 2500
$ws.Cells.Item(6, 12).Value = 2500
$ws.Cells.Item(6, 13).Value = 2500
$ws.Cells.Item(6, 16).Value = 2500
# Row 7
$ws.Cells.Item(7, 9).Value = 'Primera'
$ws.Cells.Item(7, 10).Value = 400
$ws.Cells.Item(7, 11).Value = 2000
$ws.Cells.Item(7, 12).Value = 2000
$ws.Cells.Item(7, 13).Value = 2000
$ws.Cells.Item(7, 16).Value = 2000
# Row 8
$ws.Cells.Item(8, 9).Value = 'Segunda'
$ws.Cells.Item(8, 10).Value = 300
$ws.Cells.Item(8, 11).Value = 1500
$ws.Cells.Item(8, 12).Value = 1500
$ws.Cells.Item(8, 13).Value = 1500
$ws.Cells.Item(8, 16).Value = 1500
# Row 9
$ws.Cells.Item(9, 4).Value = 44223
$ws.Cells.Item(9, 8).Value = 'Americana O Klondike'
$ws.Cells.Item(9, 9).Value = 'Tercera'
$ws.Cells.Item(9, 10).Value = 160
$ws.Cells.Item(9, 11).Value = 1000
$ws.Cells.Item(9, 12).Value = 1000
$ws.Cells.Item(9, 13).Value = 1000
$ws.Cells.Item(9, 14).Value = '$/unidad'
$ws.Cells.Item(9, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(9, 16).Value = 1000
# Row 10
$ws.Cells.Item(10, 4).Value = 44495
# Row 11
$ws.Cells.Item(11, 4).Value = 44217
$ws.Cells.Item(11, 9).Value = 'Extra'
$ws.Cells.Item(11, 10).Value = 400
$ws.Cells.Item(11, 11).Value = 2500
$ws.Cells.Item(11, 12).Value = 2500
$ws.Cells.Item(11, 13).Value = 2500
$ws.Cells.Item(11, 14).Value = '$/unidad'
$ws.Cells.Item(11, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(11, 16).Value = 2500
# Row 12
$ws.Cells.Item(12, 4).Value = 44217
$ws.Cells.Item(12, 10).Value = 280
$ws.Cells.Item(12, 11).Value = 2000
$ws.Cells.Item(12, 12).Value = 2000
$ws.Cells.Item(12, 13).Value = 2000
$ws.Cells.Item(12, 14).Value = '$/unidad'
$ws.Cells.Item(12, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(12, 16).Value = 2000
# Row 13
$ws.Cells.Item(13, 4).Value = 44194
$ws.Cells.Item(13, 10).Value = 120
$ws.Cells.Item(13, 11).Value = 3500
$ws.Cells.Item(13, 12).Value = 3500
$ws.Cells.Item(13, 13).Value = 3500
$ws.Cells.Item(13, 16).Value = 3500
# Row 14
$ws.Cells.Item(14, 4).Value = 44194
$ws.Cells.Item(14, 10).Value = 200
$ws.Cells.Item(14, 11).Value = 3000
$ws.Cells.Item(14, 12).Value = 3000
$ws.Cells.Item(14, 13).Value = 3000
$ws.Cells.Item(14, 16).Value = 3000
# Row 16
$ws.Cells.Item(16, 4).Value = 44488
$ws.Cells.Item(16, 10).Value = 150
$ws.Cells.Item(16, 11).Value = 800
$ws.Cells.Item(16, 12).Value = 800
$ws.Cells.Item(16, 13).Value = 800
$ws.Cells.Item(16, 14).Value = '$/kilo (volumen en unidades)'
$ws.Cells.Item(16, 16).Value = 800
# Row 17
$ws.Cells.Item(17, 4).Value = 44510
$ws.Cells.Item(17, 10).Value = 250
$ws.Cells.Item(17, 11).Value = 800
$ws.Cells.Item(17, 12).Value = 800
$ws.Cells.Item(17, 13).Value = 800
$ws.Cells.Item(17, 14).Value = '$/kilo (volumen en unidades)'
$ws.Cells.Item(17, 15).Value = 'Perú'
$ws.Cells.Item(17, 16).Value = 800
# Row 18
$ws.Cells.Item(18, 9).Value = 'Primera'
$ws.Cells.Item(18, 10).Value = 400
$ws.Cells.Item(18, 11).Value = 5000
$ws.Cells.Item(18, 12).Value = 5000
$ws.Cells.Item(18, 13).Value = 5000
$ws.Cells.Item(18, 16).Value = 5000
# Row 19
$ws.Cells.Item(19, 9).Value = 'Segunda'
$ws.Cells.Item(19, 10).Value = 560
$ws.Cells.Item(19, 11).Value = 3000
$ws.Cells.Item(19, 12).Value = 3000
$ws.Cells.Item(19, 13).Value = 3000
$ws.Cells.Item(19, 16).Value = 3000
# Row 20
$ws.Cells.Item(20, 4).Value = 44167
$ws.Cells.Item(20, 9).Value = 'Tercera'
$ws.Cells.Item(20, 10).Value = 450
$ws.Cells.Item(20, 11).Value = 2000
$ws.Cells.Item(20, 12).Value = 2000
$ws.Cells.Item(20, 13).Value = 2000
$ws.Cells.Item(20, 14).Value = '$/unidad'
$ws.Cells.Item(20, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(20, 16).Value = 2000
# Row 21
$ws.Cells.Item(21, 4).Value = 44491
$ws.Cells.Item(21, 10).Value = 150
$ws.Cells.Item(21, 11).Value = 800
$ws.Cells.Item(21, 12).Value = 800
$ws.Cells.Item(21, 13).Value = 800
$ws.Cells.Item(21, 14).Value = '$/kilo (volumen en unidades)'
$ws.Cells.Item(21, 16).Value = 800
# Row 22
$ws.Cells.Item(22, 4).Value = 44312
$ws.Cells.Item(22, 10).Value = 180
$ws.Cells.Item(22, 11).Value = 2500
$ws.Cells.Item(22, 12).Value = 2500
$ws.Cells.Item(22, 13).Value = 2500
$ws.Cells.Item(22, 16).Value = 2500
